# [Feat 2269] Add support of steps worksheet metadata.
# Rename the "STEPS TODO" placeholder sheet to "STEPS" and populate its
# header row with the new TC_STEP_* / TC_OWNER_* columns (replacing the
# former string-literal placeholders with the matching enum values),
# then make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# The "STEPS TODO" sheet is the 2nd sheet in the workbook.
$steps = $wb.Worksheets.Item(2)
$steps.Name = "STEPS"

# Header row - order mirrors the TC_STEPS export column layout.
$steps.Range("A1").Value = "ACTION"
$steps.Range("D1").Value = "TC_STEP_ID"
$steps.Range("E1").Value = "TC_STEP_NUM"
$steps.Range("G1").Value = "TC_STEP_ACTION"
$steps.Range("I1").Value = "TC_STEP_#_REQ"
$steps.Range("J1").Value = "TC_STEP_#_ATTACHMENT"
$steps.Range("K1").Value = "TC_STEP_CUF_<CODE>"
$steps.Range("B1").Value = "TC_OWNER_PATH"
$steps.Range("C1").Value = "TC_OWNER_ID"
$steps.Range("F1").Value = "TC_STEP_IS_CALL_STEP"
$steps.Range("H1").Value = "TC_STEP_EXPECTED_RESULT"

# Size the new columns to fit their header text.
$steps.Range("A1:K1").EntireColumn.AutoFit()

# Make STEPS the active/selected sheet and restore the previous selection.
$steps.Activate()
[void]$steps.Range("K18").Select()
